$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column values stay text (avoid numeric auto-conversion of e.g. "1.000")
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('B2').Value = 'Bitcoin'
$ws.Range('C2').Value = 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc'
$ws.Range('D2').Value = '29.418.22'
$ws.Range('E2').Value = '  +0.08%  '
$ws.Range('B3').Value = 'Ethereum'
$ws.Range('C3').Value = 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth'
$ws.Range('D3').Value = '1.848.94'
$ws.Range('E3').Value = '  +0.26%  '
$ws.Range('B4').Value = 'TetherUSD'
$ws.Range('C4').Value = 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt'
$ws.Range('D4').Value = '0.9994'
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('B5').Value = 'BNB'
$ws.Range('C5').Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range('D5').Value = '240.63'
$ws.Range('E5').Value = '  +0.73%  '
$ws.Range('B6').Value = 'XRP'
$ws.Range('C6').Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range('D6').Value = '0.6281'
$ws.Range('E6').Value = '  -0.46%  '
$ws.Range('B7').Value = 'USDC'
$ws.Range('C7').Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range('D7').Value = '1.000'
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('B8').Value = 'Dogecoin'
$ws.Range('C8').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D8').Value = '0.07718'
$ws.Range('E8').Value = '  +2.49%  '
$ws.Range('B9').Value = 'Cardano'
$ws.Range('C9').Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range('D9').Value = '0.2922'
$ws.Range('E9').Value = '  -0.30%  '
$ws.Range('B10').Value = 'Solana'
$ws.Range('C10').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('D10').Value = '25.07'
$ws.Range('E10').Value = '  +2.19%  '
$ws.Range('B11').Value = 'TRON'
$ws.Range('C11').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D11').Value = '0.07747'
$ws.Range('E11').Value = '  +0.45%  '
$ws.Range('B12').Value = 'WrappedEther'
$ws.Range('C12').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D12').Value = '1.867.13'
$ws.Range('E12').Value = '  +1.32%  '
$ws.Range('B13').Value = 'Polkadot'
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D13').Value = '5.035'
$ws.Range('E13').Value = '  +0.68%  '
$ws.Range('B14').Value = 'ShibaInu'
$ws.Range('C14').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D14').Value = '0.00001083'
$ws.Range('E14').Value = '  +3.70%  '
$ws.Range('B15').Value = 'Polygon'
$ws.Range('C15').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D15').Value = '0.6820'
$ws.Range('E15').Value = '  +0.34%  '
$ws.Range('B16').Value = 'Litecoin'
$ws.Range('C16').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D16').Value = '83.60'
$ws.Range('E16').Value = '  +0.31%  '
$ws.Range('B17').Value = 'Uniswap'
$ws.Range('C17').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D17').Value = '6.190'
$ws.Range('E17').Value = '  +0.30%  '
$ws.Range('B18').Value = 'WrappedBTC'
$ws.Range('C18').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D18').Value = '29.430.48'
$ws.Range('E18').Value = '  +0.02%  '
$ws.Range('B19').Value = 'BitcoinCash'
$ws.Range('C19').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D19').Value = '228.98'
$ws.Range('E19').Value = '  +0.20%  '
$ws.Range('B20').Value = 'Avalanche'
$ws.Range('C20').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D20').Value = '12.40'
$ws.Range('E20').Value = '  -0.10%  '
$ws.Range('B21').Value = 'Dai'
$ws.Range('C21').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D21').Value = '1.000'
$ws.Range('E21').Value = '  +0.00%  '
$ws.Range('B22').Value = 'Chainlink'
$ws.Range('C22').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D22').Value = '7.444'
$ws.Range('E22').Value = '  -0.22%  '
$ws.Range('B23').Value = 'BinanceUSD'
$ws.Range('C23').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D23').Value = '1.000'
$ws.Range('E23').Value = '  -0.04%  '
$ws.Range('B24').Value = 'Monero'
$ws.Range('C24').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D24').Value = '157.65'
$ws.Range('E24').Value = '  +0.42%  '
$ws.Range('B25').Value = 'Stellar'
$ws.Range('C25').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D25').Value = '0.1378'
$ws.Range('E25').Value = '  -1.07%  '
$ws.Range('B26').Value = 'Cosmos'
$ws.Range('C26').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D26').Value = '8.429'
$ws.Range('E26').Value = '  +0.89%  '
$ws.Range('B27').Value = 'EthereumClassic'
$ws.Range('C27').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D27').Value = '17.71'
$ws.Range('E27').Value = '  +0.70%  '
$ws.Range('B28').Value = 'Toncoin'
$ws.Range('C28').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D28').Value = '1.347'
$ws.Range('E28').Value = '  +5.10%  '
$ws.Range('B29').Value = 'PancakeSwap'
$ws.Range('C29').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D29').Value = '1.464'
$ws.Range('E29').Value = '  +0.42%  '
$ws.Range('B30').Value = 'Hedera'
$ws.Range('C30').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D30').Value = '0.05642'
$ws.Range('E30').Value = '  +0.20%  '
$ws.Range('B31').Value = 'Filecoin'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D31').Value = '4.127'
$ws.Range('E31').Value = '  +0.60%  '
$ws.Range('B32').Value = 'InternetComputer(DFINITY)'
$ws.Range('C32').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D32').Value = '4.046'
$ws.Range('E32').Value = '  +0.45%  '
$ws.Range('B33').Value = 'LidoDAOToken'
$ws.Range('C33').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D33').Value = '1.845'
$ws.Range('E33').Value = '  -0.08%  '
$ws.Range('B34').Value = 'ARBITRUM'
$ws.Range('C34').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D34').Value = '1.165'
$ws.Range('E34').Value = '  +0.77%  '
$ws.Range('B35').Value = 'ImmutableX'
$ws.Range('C35').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D35').Value = '0.7073'
$ws.Range('E35').Value = '  -0.73%  '
$ws.Range('B36').Value = 'HuobiToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D36').Value = '2.594'
$ws.Range('E36').Value = '  +0.16%  '
$ws.Range('B37').Value = 'Maker'
$ws.Range('C37').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D37').Value = '1.226.85'
$ws.Range('E37').Value = '  -1.65%  '
$ws.Range('B38').Value = 'MXToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D38').Value = '2.761'
$ws.Range('E38').Value = '  -0.21%  '
$ws.Range('B39').Value = 'VeChain'
$ws.Range('C39').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D39').Value = '0.01789'
$ws.Range('E39').Value = '  -1.22%  '
$ws.Range('B40').Value = 'FraxShare'
$ws.Range('C40').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D40').Value = '6.450'
$ws.Range('E40').Value = '  +1.22%  '
$ws.Range('B41').Value = 'TrustWalletToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D41').Value = '0.9061'
$ws.Range('E41').Value = '  +0.33%  '
$ws.Range('B42').Value = 'PaxDollar'
$ws.Range('C42').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D42').Value = '1.001'
$ws.Range('E42').Value = '  +0.06%  '
$ws.Range('B43').Value = 'RocketPoolETH'
$ws.Range('C43').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D43').Value = '2.001.66'
$ws.Range('E43').Value = '  -0.07%  '
$ws.Range('B44').Value = 'Quant'
$ws.Range('C44').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D44').Value = '101.80'
$ws.Range('E44').Value = '  +0.12%  '
$ws.Range('B45').Value = 'Aave'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D45').Value = '66.12'
$ws.Range('E45').Value = '  +0.54%  '
$ws.Range('B46').Value = 'BabyDogeCoin'
$ws.Range('C46').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D46').Value = '0.00000000122'
$ws.Range('E46').Value = '  +2.49%  '
$ws.Range('B47').Value = 'Aptos'
$ws.Range('C47').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D47').Value = '7.195'
$ws.Range('E47').Value = '  +1.37%  '
$ws.Range('B48').Value = 'TheSandbox'
$ws.Range('C48').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D48').Value = '0.4023'
$ws.Range('E48').Value = '  +0.74%  '
$ws.Range('B49').Value = 'Algorand'
$ws.Range('C49').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D49').Value = '0.1157'
$ws.Range('E49').Value = '  +3.13%  '
$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D50').Value = '9.030'
$ws.Range('E50').Value = '  +1.03%  '
$ws.Range('B51').Value = 'RenderToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D51').Value = '1.679'
$ws.Range('E51').Value = '  +0.42%  '
